$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.170.32"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "'2.310.61"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'302.55"
$ws.Range("E5").Value = "  +1.46%  "
$ws.Range("D6").Value = "'100.05"
$ws.Range("E6").Value = "  +5.70%  "
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.508"
$ws.Range("E9").Value = "  +3.50%  "
$ws.Range("D10").Value = "'34.44"
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("D11").Value = "'0.0797"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'0.117"
$ws.Range("E12").Value = "  +3.78%  "
$ws.Range("D13").Value = "'17.91"
$ws.Range("E13").Value = "  +14.45%  "
$ws.Range("D14").Value = "'6.83"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "'2.669.67"
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.819"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.267.57"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'43.117.83"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").Value = "'12.66"
$ws.Range("E19").Value = "  +11.41%  "
$ws.Range("D20").Value = "'0.0₃0907"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("D22").Value = "'67.83"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "'237.36"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").Value = "'2.20"
$ws.Range("E24").Value = "  +13.56%  "
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'24.78"
$ws.Range("E27").Value = "  +4.21%  "
$ws.Range("E28").Value = "  -5.06%  "
$ws.Range("D29").Value = "'168.54"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("D30").Value = "'34.20"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "'9.17"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'5.04"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").Value = "'4.61"
$ws.Range("E34").Value = "  +3.37%  "
$ws.Range("E35").Value = "  +4.65%  "
$ws.Range("D36").Value = "'17.14"
$ws.Range("E36").Value = "  +6.93%  "
$ws.Range("D37").Value = "'0.0692"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("E38").Value = "  +3.87%  "
$ws.Range("D39").Value = "'1.80"
$ws.Range("E39").Value = "  +5.09%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'2.000.72"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "'10.14"
$ws.Range("E45").Value = "  +6.16%  "
$ws.Range("D46").Value = "'17.64"
$ws.Range("D47").Value = "'2.86"
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("D48").Value = "'56.02"
$ws.Range("E48").Value = "  +7.92%  "
$ws.Range("D49").Value = "'2.537.66"
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").Value = "'4.57"
$ws.Range("E51").Value = "  +1.40%  "
